# account_range.xlsx - add the missing F-3 batch row (03200501200101-03200501200152)
# so the QR generator has real account data for that batch instead of
# skipping/erroring on an empty lookup.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: S.N becomes 2
$ws.Range("A3").Value = 2

# Row 3: Start column (E3) becomes text "03200501200100" (leading zero preserved)
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "03200501200100"

# Row 4: new account range row
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "F-3"
$ws.Range("D4").Value = 3200501200101
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "03200501200152"
$ws.Range("F4").Value = "F-02-03"
$ws.Range("C4").Value = "F- 33"
$ws.Range("G4").Value = 12

# Update selection to match target state
$ws.Range("C4").Select() | Out-Null
